# The base bus-load value lives in cell B5 of the "buses" sheet. Every
# other cell in column B (B6:B104) holds a formula that chains back to
# it (B6 = B5, B7 = B6, ... B104 = B103), so changing B5 and letting
# Excel recalculate ripples the new value down the whole column - exactly
# like a user editing the cell in the live workbook.
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "buses") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Range("B5").Value = -0.01

$excel.CalculateFullRebuild()
